# Sweden Superettan workbook update ("Atualização de bases das ligas, do dia: 2024-01-29 às 18-52")
#
# The upstream scrape re-sorted a handful of match rows that share an
# identical Date value (column E). Excel/pandas's sort is not guaranteed
# stable for ties, so several adjacent rows (same date) swapped places:
# every column except A (the permanent row sequence number), C, D and E
# (league name / date, identical for the whole block) moved together as a
# unit from one row to another inside each small block of tied rows.
#
# Implementation: for every affected block of 2 (or 3) adjacent rows we
# snapshot the B:AC range of each row *before* writing anything (so the
# write of one row never clobbers data still needed for another row in
# the same block), then write each row's new B:AC content from the
# snapshot of the row whose data it should now hold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: list of adjacent row numbers that form one rotation "block",
# followed by the permutation describing, for every row in the block, which
# row's *original* data it should end up holding.
$blocks = @(
    @{ Rows = @(344, 345);      Source = @{ 344 = 345; 345 = 344 } },
    @{ Rows = @(360, 361);      Source = @{ 360 = 361; 361 = 360 } },
    @{ Rows = @(370, 371, 372); Source = @{ 370 = 371; 371 = 372; 372 = 370 } },
    @{ Rows = @(392, 393);      Source = @{ 392 = 393; 393 = 392 } },
    @{ Rows = @(416, 417);      Source = @{ 416 = 417; 417 = 416 } },
    @{ Rows = @(432, 433);      Source = @{ 432 = 433; 433 = 432 } },
    @{ Rows = @(443, 444);      Source = @{ 443 = 444; 444 = 443 } },
    @{ Rows = @(448, 449);      Source = @{ 448 = 449; 449 = 448 } },
    @{ Rows = @(454, 455, 456); Source = @{ 454 = 456; 455 = 454; 456 = 455 } },
    @{ Rows = @(470, 471);      Source = @{ 470 = 471; 471 = 470 } },
    @{ Rows = @(484, 485);      Source = @{ 484 = 485; 485 = 484 } },
    @{ Rows = @(706, 707);      Source = @{ 706 = 707; 707 = 706 } }
)

foreach ($block in $blocks) {
    # 1) Snapshot every row's current B:AC values first.
    $snapshots = @{}
    foreach ($row in $block.Rows) {
        $snapshots[$row] = $ws.Range("B$($row):AC$($row)").Value()
    }

    # 2) Now write each row's new content from the recorded snapshots.
    foreach ($row in $block.Rows) {
        $srcRow = $block.Source[$row]
        $ws.Range("B$($row):AC$($row)").Value = $snapshots[$srcRow]
    }
}
